# Omaha_Cal_Info_GA05MOAS-GL002_00001.xlsx
# - Asset_Cal_Info!F2 (CC_scattering_angle): 117 -> 140
# - Asset_Cal_Info!F4 (CC_angular_resolution): 1.08 -> 1.13
# - Active sheet moves from Moorings to Asset_Cal_Info, with the
#   selection left on F4 (matches the saved workbook/sheet view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asset_Cal_Info")

$ws.Range("F2").Value = 140
$ws.Range("F4").Value = 1.13

$ws.Activate()
$ws.Range("F4").Select()
